$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.127881588408715
$ws.Range("C2").Value = 0.002777888934908601
$ws.Range("D2").Value = 3.900430680208489
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 12.69132264350109
